$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was D2=44881, L2=Segunda, M2=100, N2=11250, O2=11250, P2=11250, S2=11250)
$ws.Range("D2").Value = 44923
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 7500
$ws.Range("O2").Value = 8000
$ws.Range("P2").Value = 7625
$ws.Range("S2").Value = 7625

# Row 3 (was D3=44874, M3=200, P3=7750, S3=7750)
$ws.Range("D3").Value = 44923
$ws.Range("M3").Value = 80
$ws.Range("P3").Value = 7625
$ws.Range("S3").Value = 7625

# Row 4 (was D4=44923, L4=Primera, M4=80, N4=7500, O4=8000, P4=7625, S4=7625)
$ws.Range("D4").Value = 44881
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 11250
$ws.Range("O4").Value = 11250
$ws.Range("P4").Value = 11250
$ws.Range("S4").Value = 11250

# Row 5 (was D5=44923, M5=80, P5=7625, S5=7625)
$ws.Range("D5").Value = 44874
$ws.Range("M5").Value = 200
$ws.Range("P5").Value = 7750
$ws.Range("S5").Value = 7750
